# Adds "manual_features" / "num_manual_features" columns to Table1 on the
# "toybox_0_7_5_infer" sheet, populates the manual-feature annotations that
# were added for a handful of rows, and expands the table/columns to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("toybox_0_7_5_infer")

# --- 1. Grow the table by two columns (line/filename/.../classification + 2) ---
$lo = $ws.ListObjects.Item("Table1")
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null

# Name the two new header cells (renaming the table header cell also renames
# the corresponding ListColumn in this engine).
$ws.Range("I1").Value = "manual_features"
$ws.Range("J1").Value = "num_manual_features"

# --- 2. A handful of data rows were re-sorted (by filename) as part of this
# edit, which shuffled a few rows among themselves; reproduce the rows whose
# A:H content changed position, then stamp the new manual_features /
# num_manual_features values (I/J) for every row that has them ---

# Row 77 <- old row 112 (mountpoint.c / MEMORY_LEAK)
$ws.Cells.Item(77, 1).Value = 53
$ws.Cells.Item(77, 2).Value = "toys/other/mountpoint.c"
$ws.Cells.Item(77, 3).Value = "infer"
$ws.Cells.Item(77, 4).Value = "toybox_0_7_5"
$ws.Cells.Item(77, 5).Value = "MEMORY_LEAK"
$ws.Cells.Item(77, 6).Value = 258
$ws.Cells.Item(77, 7).Value = $true
$ws.Cells.Item(77, 8).Value = $true

# Row 78 <- old row 77 (oneit.c / RESOURCE_LEAK, line 72)
$ws.Cells.Item(78, 1).Value = 72
$ws.Cells.Item(78, 2).Value = "toys/other/oneit.c"
$ws.Cells.Item(78, 3).Value = "infer"
$ws.Cells.Item(78, 4).Value = "toybox_0_7_5"
$ws.Cells.Item(78, 5).Value = "RESOURCE_LEAK"
$ws.Cells.Item(78, 6).Value = 516
$ws.Cells.Item(78, 7).Value = $true
$ws.Cells.Item(78, 8).Value = $true

# Row 95 <- old row 78 (oneit.c / RESOURCE_LEAK, line 73)
$ws.Cells.Item(95, 1).Value = 73
$ws.Cells.Item(95, 2).Value = "toys/other/oneit.c"
$ws.Cells.Item(95, 3).Value = "infer"
$ws.Cells.Item(95, 4).Value = "toybox_0_7_5"
$ws.Cells.Item(95, 5).Value = "RESOURCE_LEAK"
$ws.Cells.Item(95, 6).Value = 516
$ws.Cells.Item(95, 7).Value = $true
$ws.Cells.Item(95, 8).Value = $true

# Row 102 <- old row 95 (uptime.c / NULL_DEREFERENCE)
$ws.Cells.Item(102, 1).Value = 54
$ws.Cells.Item(102, 2).Value = "toys/other/uptime.c"
$ws.Cells.Item(102, 3).Value = "infer"
$ws.Cells.Item(102, 4).Value = "toybox_0_7_5"
$ws.Cells.Item(102, 5).Value = "NULL_DEREFERENCE"
$ws.Cells.Item(102, 6).Value = 501
$ws.Cells.Item(102, 7).Value = $true
$ws.Cells.Item(102, 8).Value = $true

# Row 103 <- old row 102 (nohup.c / RESOURCE_LEAK, line 27)
$ws.Cells.Item(103, 1).Value = 27
$ws.Cells.Item(103, 2).Value = "toys/posix/nohup.c"
$ws.Cells.Item(103, 3).Value = "infer"
$ws.Cells.Item(103, 4).Value = "toybox_0_7_5"
$ws.Cells.Item(103, 5).Value = "RESOURCE_LEAK"
$ws.Cells.Item(103, 6).Value = 485
$ws.Cells.Item(103, 7).Value = $true
$ws.Cells.Item(103, 8).Value = $true

# Row 112 <- old row 103 (nohup.c / RESOURCE_LEAK, line 39)
$ws.Cells.Item(112, 1).Value = 39
$ws.Cells.Item(112, 2).Value = "toys/posix/nohup.c"
$ws.Cells.Item(112, 3).Value = "infer"
$ws.Cells.Item(112, 4).Value = "toybox_0_7_5"
$ws.Cells.Item(112, 5).Value = "RESOURCE_LEAK"
$ws.Cells.Item(112, 6).Value = 485
$ws.Cells.Item(112, 7).Value = $true
$ws.Cells.Item(112, 8).Value = $true

# --- 3. Stamp manual_features / num_manual_features (I/J) on every row that
# carries one, including rows whose A:H content did not move ---
$ws.Cells.Item(77, 9).Value = "CONFIG_MOUNTPOINT, -CONFIG_TOYBOX_FREE"
$ws.Cells.Item(77, 10).Value = 2

$ws.Cells.Item(78, 9).Value = "CONFIG_ONEIT"
$ws.Cells.Item(78, 10).Value = 1

$ws.Cells.Item(79, 9).Value = "CONFIG_ONEIT"

$ws.Cells.Item(95, 9).Value = "CONFIG_ONEIT"
$ws.Cells.Item(95, 10).Value = 1

$ws.Cells.Item(96, 9).Value = "CONFIG_ONEIT"

$ws.Cells.Item(102, 9).Value = "CONFIG_UPTIME"
$ws.Cells.Item(102, 10).Value = 1

$ws.Cells.Item(103, 9).Value = "CONFIG_NOHUP"
$ws.Cells.Item(103, 10).Value = 1

$ws.Cells.Item(104, 9).Value = "CONFIG_NOHUP"

$ws.Cells.Item(112, 9).Value = "CONFIG_NOHUP"
$ws.Cells.Item(112, 10).Value = 1

$ws.Cells.Item(113, 9).Value = "CONFIG_NOHUP"

# --- 4. Column layout: widen filename, hide tool/target/type, widen the new
# manual_features column ---
$ws.Columns.Item(2).ColumnWidth = 22.6640625
$ws.Range("C1:E1").EntireColumn.Hidden = $true
$ws.Columns.Item(9).ColumnWidth = 42.6640625

# --- 5. View state: the edited sheet was scrolled right and K2 selected ---
$ws.Range("K2").Select()
